$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1.632
$ws.Range("D2").Value = 1.662
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 20
$ws.Range("G2").Value = 29
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 1103.7
$ws.Range("J2").Value = -99.84941560206579

$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 476
$ws.Range("D3").Value = 485
$ws.Range("E3").Value = 5
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 7
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 541
$ws.Range("J3").Value = -10.35120147874307

$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = 21
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 12
$ws.Range("J4").Value = 75

$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 442
$ws.Range("D5").Value = 450
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 219
$ws.Range("J5").Value = 105.4794520547945

$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 68
$ws.Range("D6").Value = 68
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 92
$ws.Range("J6").Value = -26.08695652173914

$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 63
$ws.Range("D7").Value = 70
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 196
$ws.Range("J7").Value = -64.28571428571428

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 390
$ws.Range("D8").Value = 417
$ws.Range("E8").Value = 27
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = 15
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 23
$ws.Range("J8").Value = 1713.04347826087

$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 74
$ws.Range("D9").Value = 75
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 58
$ws.Range("J9").Value = 29.31034482758621

$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 302
$ws.Range("D10").Value = 356
$ws.Range("E10").Value = 45
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 7
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 915
$ws.Range("J10").Value = -61.09289617486338

$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 391
$ws.Range("D11").Value = 396
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 7
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 372
$ws.Range("J11").Value = 6.451612903225801

$ws.Range("B12").Value = 23
$ws.Range("C12").Value = 387
$ws.Range("D12").Value = 1129
$ws.Range("E12").Value = 462
$ws.Range("F12").Value = 18
$ws.Range("G12").Value = 7
$ws.Range("H12").Value = 239
$ws.Range("I12").Value = 1019.2
$ws.Range("J12").Value = 10.77315541601256

$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 26
$ws.Range("D13").Value = 44
$ws.Range("E13").Value = 17
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 766
$ws.Range("J13").Value = -94.25587467362925

$ws.Range("B14").Value = 11
$ws.Range("C14").Value = 376
$ws.Range("D14").Value = 793
$ws.Range("E14").Value = 339
$ws.Range("F14").Value = 10
$ws.Range("G14").Value = 6
$ws.Range("H14").Value = 58
$ws.Range("I14").Value = 1157
$ws.Range("J14").Value = -31.46067415730337

$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 254
$ws.Range("D15").Value = 278
$ws.Range("E15").Value = 24
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 309
$ws.Range("J15").Value = -10.03236245954693

$ws.Range("C17").Value = 70
$ws.Range("D17").Value = 95
$ws.Range("E17").Value = 25
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 1
$ws.Range("I17").Value = 191
$ws.Range("J17").Value = -50.26178010471204

$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 2
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 0

$ws.Range("C19").Value = 23
$ws.Range("D19").Value = 23
$ws.Range("I19").Value = 22
$ws.Range("J19").Value = 4.545454545454541

$ws.Range("C20").Value = 58
$ws.Range("D20").Value = 64
$ws.Range("E20").Value = 6
$ws.Range("I20").Value = 108
$ws.Range("J20").Value = -40.74074074074075
